# Apply cryptos list update (Sat Oct  5 07:29:10 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.203.10"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3
$ws.Range("D3").Value = "2.422.12"
$ws.Range("E3").Value = "  +1.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("E8").Value = "  +1.90%  "

# Row 9
$ws.Range("D9").Value = "2.421.63"
$ws.Range("E9").Value = "  +1.76%  "

# Row 10
$ws.Range("E10").Value = "  +1.29%  "

# Row 11
$ws.Range("E11").Value = "  -1.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "

# Row 13
$ws.Range("E13").Value = "  +1.19%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.77%  "

# Row 15
$ws.Range("E15").Value = "  +5.63%  "

# Row 16
$ws.Range("E16").Value = "  +2.16%  "

# Row 17
$ws.Range("D17").Value = "62.035.80"
$ws.Range("E17").Value = "  +1.32%  "

# Row 18
$ws.Range("D18").Value = "2.420.36"
$ws.Range("E18").Value = "  +2.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21
$ws.Range("E21").Value = "  +0.83%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.89%  "

# Row 25
$ws.Range("E25").Value = "  +0.54%  "

# Row 26
$ws.Range("E26").Value = "  +5.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "585.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.03%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.80%  "

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.542.88"
$ws.Range("E29").Value = "  +2.16%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0944"
$ws.Range("E30").Value = "  +5.08%  "

# Row 31
$ws.Range("E31").Value = "  +1.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.150"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.61%  "

# Row 34
$ws.Range("E34").Value = "  +3.29%  "

# Row 35
$ws.Range("E35").Value = "  +2.01%  "

# Row 36
$ws.Range("E36").Value = "  +4.79%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("E38").Value = "  +1.86%  "

# Row 39
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.37%  "

# Row 40
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.385"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.62%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.95%  "

# Row 42
$ws.Range("E42").Value = "  -4.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

# Row 44
$ws.Range("E44").Value = "  +9.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46
$ws.Range("E46").Value = "  +1.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0540"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.96%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.593"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.54%  "

# Row 50
$ws.Range("E50").Value = "  +2.14%  "

# Row 51
$ws.Range("E51").Value = "  +2.26%  "

